$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.83803220311825
$ws.Range("C2").Value = 9.178809736610043
$ws.Range("D2").Value = 8.461251432989892
$ws.Range("F2").Value = 38.84545948736759
$ws.Range("G2").Value = 3.697694011909839
$ws.Range("J2").Value = 10.99748134173379
$ws.Range("K2").Value = 10.81443990301441
$ws.Range("L2").Value = 11.37542160696269
$ws.Range("M2").Value = 16.16178540163461
$ws.Range("N2").Value = 21.99220027066415
$ws.Range("O2").Value = 29.76247574023632
$ws.Range("B3").Value = 14.68204625040448
$ws.Range("C3").Value = 9.167035933157711
$ws.Range("D3").Value = 8.454046083180176
$ws.Range("F3").Value = 38.92682386148385
$ws.Range("G3").Value = 3.699464746767345
$ws.Range("J3").Value = 11.01697295416666
$ws.Range("K3").Value = 10.69719650617781
$ws.Range("L3").Value = 11.38400891367576
$ws.Range("M3").Value = 16.14516534871521
$ws.Range("N3").Value = 22.05247610616572
$ws.Range("O3").Value = 29.84027972293675
$ws.Range("B4").Value = 14.58819134287216
$ws.Range("C4").Value = 9.15983613829362
$ws.Range("D4").Value = 8.450710100676599
$ws.Range("F4").Value = 38.98362720418062
$ws.Range("G4").Value = 3.700610772019758
$ws.Range("J4").Value = 11.02963073672081
$ws.Range("K4").Value = 10.62636396506587
$ws.Range("L4").Value = 11.39049366030994
$ws.Range("M4").Value = 16.13706023727897
$ws.Range("N4").Value = 22.09123746158041
$ws.Range("O4").Value = 29.89282435771996
$ws.Range("B5").Value = 14.55046835658555
$ws.Range("C5").Value = 9.156910093184417
$ws.Range("D5").Value = 8.449625881899205
$ws.Range("F5").Value = 39.00849528854188
$ws.Range("G5").Value = 3.701092614149061
$ws.Range("J5").Value = 11.0349628103925
$ws.Range("K5").Value = 10.59781978800441
$ws.Range("L5").Value = 11.39344174558853
$ws.Range("M5").Value = 16.1342886416389
$ws.Range("N5").Value = 22.10747477985989
$ws.Range("O5").Value = 29.91543637915532
$ws.Range("B6").Value = 14.54423727778763
$ws.Range("C6").Value = 9.156424720373154
$ws.Range("D6").Value = 8.449462515353121
$ws.Range("F6").Value = 39.01272847981353
$ws.Range("G6").Value = 3.701173520540904
$ws.Range("J6").Value = 11.03585871518624
$ws.Range("K6").Value = 10.59310028569967
$ws.Range("L6").Value = 11.39394974410354
$ws.Range("M6").Value = 16.1338606019787
$ws.Range("N6").Value = 22.11019769195258
$ws.Range("O6").Value = 29.91926352956344
$ws.Range("B7").Value = 14.5876804243409
$ws.Range("C7").Value = 9.159796643799949
$ws.Range("D7").Value = 8.450694362127374
$ws.Range("F7").Value = 38.98395562013636
$ws.Range("G7").Value = 3.700617210214997
$ws.Range("J7").Value = 11.02970194211023
$ws.Range("K7").Value = 10.62597767058483
$ws.Range("L7").Value = 11.39053218136124
$ws.Range("M7").Value = 16.13702070313589
$ws.Range("N7").Value = 22.09145465343002
$ws.Range("O7").Value = 29.89312445465663
$ws.Range("B8").Value = 14.78387413530912
$ws.Range("C8").Value = 9.174744359973143
$ws.Range("D8").Value = 8.458542180820848
$ws.Range("F8").Value = 38.87209261874231
$ws.Range("G8").Value = 3.698292387250518
$ws.Range("J8").Value = 11.00405914475688
$ws.Range("K8").Value = 10.77379281172916
$ws.Range("L8").Value = 11.3781313704873
$ws.Range("M8").Value = 16.15562108502388
$ws.Range("N8").Value = 22.01262056514135
$ws.Range("O8").Value = 29.7883117175542
$ws.Range("B9").Value = 15.18203165550558
$ws.Range("C9").Value = 9.204266719239303
$ws.Range("D9").Value = 8.482491136826566
$ws.Range("F9").Value = 38.70708421662306
$ws.Range("G9").Value = 3.694197807939072
$ws.Range("J9").Value = 10.95922668653717
$ws.Range("K9").Value = 11.07148056817659
$ws.Range("L9").Value = 11.36339925314339
$ws.Range("M9").Value = 16.20860164337799
$ws.Range("N9").Value = 21.87186801080452
$ws.Range("O9").Value = 29.62066318914658
$ws.Range("B10").Value = 15.48026673501316
$ws.Range("C10").Value = 9.226055032320621
$ws.Range("D10").Value = 8.505198725904632
$ws.Range("F10").Value = 38.61903405321554
$ws.Range("G10").Value = 3.69146974998055
$ws.Range("J10").Value = 10.92958372274664
$ws.Range("K10").Value = 11.29314349717978
$ws.Range("L10").Value = 11.35837568286464
$ws.Range("M10").Value = 16.2573685730255
$ws.Range("N10").Value = 21.7768115726333
$ws.Range("O10").Value = 29.52061521480807
$ws.Range("B11").Value = 15.61665238889362
$ws.Range("C11").Value = 9.235982039859822
$ws.Range("D11").Value = 8.516614121981689
$ws.Range("F11").Value = 38.58618771919673
$ws.Range("G11").Value = 3.690288923882515
$ws.Range("J11").Value = 10.91680783589521
$ws.Range("K11").Value = 11.3942388087466
$ws.Range("L11").Value = 11.35733967056678
$ws.Range("M11").Value = 16.28163915344726
$ws.Range("N11").Value = 21.73536425728195
$ws.Range("O11").Value = 29.48012613818382
$ws.Range("B12").Value = 15.66835720768004
$ws.Range("C12").Value = 9.239742763131495
$ws.Range("D12").Value = 8.521090612212888
$ws.Range("F12").Value = 38.57478609813399
$ws.Range("G12").Value = 3.689850383503702
$ws.Range("J12").Value = 10.91207141056922
$ws.Range("K12").Value = 11.43252638307825
$ws.Range("L12").Value = 11.35712605636413
$ws.Range("M12").Value = 16.29112495555424
$ws.Range("N12").Value = 21.71992595139318
$ws.Range("O12").Value = 29.46551658226513
$ws.Range("B13").Value = 15.65721974998985
$ws.Range("C13").Value = 9.238932763925121
$ws.Range("D13").Value = 8.520119722041885
$ws.Range("F13").Value = 38.57719554161381
$ws.Range("G13").Value = 3.689944448613859
$ws.Range("J13").Value = 10.91308697553613
$ws.Range("K13").Value = 11.4242807574524
$ws.Range("L13").Value = 11.35716412964939
$ws.Range("M13").Value = 16.28906897672852
$ws.Range("N13").Value = 21.72323945929354
$ws.Range("O13").Value = 29.4686308595633
$ws.Range("B14").Value = 15.62090525560575
$ws.Range("C14").Value = 9.236291408791491
$ws.Range("D14").Value = 8.51697933865564
$ws.Range("F14").Value = 38.58522892723093
$ws.Range("G14").Value = 3.690252672522576
$ws.Range("J14").Value = 10.91641613491311
$ws.Range("K14").Value = 11.3973888343934
$ws.Range("L14").Value = 11.35731852049607
$ws.Range("M14").Value = 16.28241367181693
$ws.Range("N14").Value = 21.73408899648404
$ws.Range("O14").Value = 29.47890971341502
$ws.Range("B15").Value = 15.59866788880767
$ws.Range("C15").Value = 9.234673689491334
$ws.Range("D15").Value = 8.515075710337184
$ws.Range("F15").Value = 38.59028459690098
$ws.Range("G15").Value = 3.690442588996514
$ws.Range("J15").Value = 10.91846855084711
$ws.Range("K15").Value = 11.38091645090675
$ws.Range("L15").Value = 11.35743633289454
$ws.Range("M15").Value = 16.27837538012527
$ws.Range("N15").Value = 21.74076807296758
$ws.Range("O15").Value = 29.48529994707792
$ws.Range("B16").Value = 15.47136421708029
$ws.Range("C16").Value = 9.225406534884383
$ws.Range("D16").Value = 8.504474347204852
$ws.Range("F16").Value = 38.62132569541841
$ws.Range("G16").Value = 3.69154812721673
$ws.Range("J16").Value = 10.93043288603938
$ws.Range("K16").Value = 11.28653912825025
$ws.Range("L16").Value = 11.35846844745866
$ws.Range("M16").Value = 16.25582399753683
$ws.Range("N16").Value = 21.77955626190569
$ws.Range("O16").Value = 29.52336239132929
$ws.Range("B17").Value = 15.39341930223304
$ws.Range("C17").Value = 9.219724930698501
$ws.Range("D17").Value = 8.498247146762546
$ws.Range("F17").Value = 38.64221470678595
$ws.Range("G17").Value = 3.692241723558773
$ws.Range("J17").Value = 10.93795388082306
$ws.Range("K17").Value = 11.22868503210945
$ws.Range("L17").Value = 11.35942092339791
$ws.Range("M17").Value = 16.24252035758215
$ws.Range("N17").Value = 21.80381032405234
$ws.Range("O17").Value = 29.5479993814527
$ws.Range("B18").Value = 15.34865757791801
$ws.Range("C18").Value = 9.216458509576745
$ws.Range("D18").Value = 8.494767757675364
$ws.Range("F18").Value = 38.65490799145254
$ws.Range("G18").Value = 3.692646328800838
$ws.Range("J18").Value = 10.94234650202256
$ws.Range("K18").Value = 11.19543524385169
$ws.Range("L18").Value = 11.36008640532165
$ws.Range("M18").Value = 16.23506513569254
$ws.Range("N18").Value = 21.81792957727548
$ws.Range("O18").Value = 29.56264271688942
$ws.Range("B19").Value = 15.33351543892273
$ws.Range("C19").Value = 9.215352830119246
$ws.Range("D19").Value = 8.493607342485877
$ws.Range("F19").Value = 38.65932224246573
$ws.Range("G19").Value = 3.692784295750993
$ws.Range("J19").Value = 10.94384524384093
$ws.Range("K19").Value = 11.18418295839748
$ws.Range("L19").Value = 11.36033195914484
$ws.Range("M19").Value = 16.2325748521985
$ws.Range("N19").Value = 21.82273916717303
$ws.Range("O19").Value = 29.56768188614248
$ws.Range("B20").Value = 15.4017097380445
$ws.Range("C20").Value = 9.220329597794615
$ws.Range("D20").Value = 8.49889946878454
$ws.Range("F20").Value = 38.63992081465285
$ws.Range("G20").Value = 3.692167302832662
$ws.Range("J20").Value = 10.93714635329401
$ws.Range("K20").Value = 11.23484120762705
$ws.Range("L20").Value = 11.35930736174705
$ws.Range("M20").Value = 16.24391623215448
$ws.Range("N20").Value = 21.80121095940973
$ws.Range("O20").Value = 29.54532779631175
$ws.Range("B21").Value = 15.63157046665262
$ws.Range("C21").Value = 9.237067200993053
$ws.Range("D21").Value = 8.517897592819621
$ws.Range("F21").Value = 38.58284119586617
$ws.Range("G21").Value = 3.690161906265033
$ws.Range("J21").Value = 10.91543552793518
$ws.Range("K21").Value = 11.40528776444821
$ws.Range("L21").Value = 11.35726832996726
$ws.Range("M21").Value = 16.28436052893295
$ws.Range("N21").Value = 21.73089526136836
$ws.Range("O21").Value = 29.47587094574138
$ws.Range("B22").Value = 15.78211869481182
$ws.Range("C22").Value = 9.248015108771401
$ws.Range("D22").Value = 8.531208819291281
$ws.Range("F22").Value = 38.55157812382301
$ws.Range("G22").Value = 3.688901447046403
$ws.Range("J22").Value = 10.90183783309078
$ws.Range("K22").Value = 11.51669882368226
$ws.Range("L22").Value = 11.35697684386547
$ws.Range("M22").Value = 16.31251058916561
$ws.Range("N22").Value = 21.68643673960608
$ws.Range("O22").Value = 29.43468993450623
$ws.Range("B23").Value = 15.70175346765545
$ws.Range("C23").Value = 9.242171399360778
$ws.Range("D23").Value = 8.524023304878762
$ws.Range("F23").Value = 38.56771104993894
$ws.Range("G23").Value = 3.6895695995724
$ws.Range("J23").Value = 10.90904118158692
$ws.Range("K23").Value = 11.45724587988976
$ws.Range("L23").Value = 11.35703747302656
$ws.Range("M23").Value = 16.29733092703744
$ws.Range("N23").Value = 21.71002849797889
$ws.Range("O23").Value = 29.45628339497241
$ws.Range("B24").Value = 15.39796147334337
$ws.Range("C24").Value = 9.220056227731231
$ws.Range("D24").Value = 8.498604240126502
$ws.Range("F24").Value = 38.64095575294978
$ws.Range("G24").Value = 3.692200930225397
$ws.Range("J24").Value = 10.93751122246211
$ws.Range("K24").Value = 11.23205796299159
$ws.Range("L24").Value = 11.35935833568605
$ws.Range("M24").Value = 16.24328455503621
$ws.Range("N24").Value = 21.80238558611402
$ws.Range("O24").Value = 29.54653412729893
$ws.Range("B25").Value = 15.07313694823101
$ws.Range("C25").Value = 9.196260704732861
$ws.Range("D25").Value = 8.475106393181834
$ws.Range("F25").Value = 38.74589862952249
$ws.Range("G25").Value = 3.695256082440443
$ws.Range("J25").Value = 10.97077427154539
$ws.Range("K25").Value = 10.9902980410689
$ws.Range("L25").Value = 11.36636278589563
$ws.Range("M25").Value = 16.28906897672852
$ws.Range("N25").Value = 21.72323945929354
$ws.Range("O25").Value = 29.4686308595633
